# Revisão da classificação feita
$wb = $excel.ActiveWorkbook
$wsTreino = $wb.Worksheets.Item("Treinamento")
$wsTeste  = $wb.Worksheets.Item("Teste")

# Fix the stray "-" placeholder in the Teste classification column to a proper 0
$wsTeste.Range("B191").Value = 0

# Review pass: flip several classification values from 0 to 1
$rowsToFlip = @(16,19,26,40,42,43,45,49,82,102,116,142,145,146,153,246,253,263,264,266,291)
foreach ($r in $rowsToFlip) {
    $wsTeste.Cells.Item($r, 2).Value = 1
}

# Make "Teste" the active sheet/tab (was "Treinamento" before)
$wsTeste.Activate()
$wb.Application.ActiveWindow.ActiveCell.Select() | Out-Null
